$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated sval data (filtered save games), rows 2-11, columns B-G.
# G (sum) = B + C + D + E (F = Win indicator, not included in sum).
$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    3  = @(0.0001488876196638067, 0.002777888934908601, 0.1575252929769615, 0.496779210170732, 1, 0.6572312797022659)
    4  = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1, 1.271902929317955)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    6  = @(0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1, 2.290389397800092)
    7  = @(1.459612070389937, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 2.42670696938877)
    8  = @(0.3048080303191223, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.626907116734944)
    9  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 0, 11.945164432584)
    10 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 1, 12.59312877619104)
    11 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
